$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetime1" field text (10/18/23 -> 10/19/23) on the
#    slide master and every slide layout (the footer "Date Placeholder").
# ---------------------------------------------------------------------------
$m = $p.SlideMaster

$masterShapeCount = $m.Shapes.Count
$msi = 1
while ($msi -le $masterShapeCount) {
    $msh = $m.Shapes.Item($msi)
    if ($msh.HasTextFrame) {
        if ($msh.TextFrame.TextRange.Text -eq "10/18/23") {
            $msh.TextFrame.TextRange.Text = "10/19/23"
        }
    }
    $msi = $msi + 1
}

$layouts = $m.CustomLayouts
$layoutCount = $layouts.Count
$li = 1
while ($li -le $layoutCount) {
    $lay = $layouts.Item($li)
    $shCount = $lay.Shapes.Count
    $si = 1
    while ($si -le $shCount) {
        $sh = $lay.Shapes.Item($si)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "10/18/23") {
                $sh.TextFrame.TextRange.Text = "10/19/23"
            }
        }
        $si = $si + 1
    }
    $li = $li + 1
}

# ---------------------------------------------------------------------------
# 2) Slide 4 ("... areas:" slide): the "Bias" / "Fairness" two-line label
#    becomes "Bias &" / "Fairness".
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$biasShape = $s4.Shapes.Item(4)
$biasRun = $biasShape.TextFrame.TextRange.Characters(1, 4)
$biasRun.Text = "Bias &"

# ---------------------------------------------------------------------------
# 3) Slide 6 ("Important dates" agenda slide).
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# 3a) "Deadline to rank topics" textbox shifts slightly to the left.
$deadlineShape = $s6.Shapes.Item(15)
$deadlineShape.Left = 229.59181213378906

# 3b) "Paper assignment / Presentation dates" textbox shifts slightly to the
#     left and its first line gains a trailing " &".
$paperShape = $s6.Shapes.Item(20)
$paperShape.Left = 325.9804992675781
$paperRun = $paperShape.TextFrame.TextRange.Characters(1, 17)
$paperRun.Text = "Paper assignment &"

# 3c) "Test presentations" becomes "Practice presentations".
$testShape = $s6.Shapes.Item(23)
$testShape.TextFrame.TextRange.Text = "Practice presentations"

# 3d) The connector sitting under the brace shifts slightly to the left too.
$connectorShape = $s6.Shapes.Item(27)
$connectorShape.Left = 269.2812805175781
